$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bang luong")

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Bùi Anh Dũng"
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 2025
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 500000
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Bùi Anh Dũng"
$ws.Cells.Item(6, 3).Value = 11
$ws.Cells.Item(6, 4).Value = 2025
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 500000
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
